$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5027148723602295
$ws.Range("B1").Value = 0.633942186832428
$ws.Range("C1").Value = 0.9215063452720642
$ws.Range("D1").Value = 3.652648448944092
$ws.Range("E1").Value = 5.622573852539062
